$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-Format($srcRef, $destRef) {
    $ws.Range($srcRef).Copy()
    $ws.Range($destRef).PasteSpecial(-4122)
}

$ws.Range("M70").Value = "Merovingok"
$ws.Hyperlinks.Add($ws.Range("AG70"), "https://hu.wikipedia.org/wiki/N%C3%A9metorsz%C3%A1g_t%C3%B6rt%C3%A9nelme") | Out-Null
Copy-Format "AG5" "AG70"
$ws.Range("I72").Value = "751-911"
Copy-Format "I7" "I72"
$ws.Range("M72").Value = "Karolingok"
Copy-Format "M7" "M72"
$ws.Range("X70").Value = "Merovingok.txt"
Copy-Format "X7" "X70"
$ws.Range("I74").Value = "919-1024"
Copy-Format "I7" "I74"
$ws.Range("M74").Value = "Az Ottók uralma"
Copy-Format "M7" "M74"
$ws.Range("X72").Value = "Karolingok.txt"
Copy-Format "I7" "X72"
$ws.Range("X74").Value = "Ottók.txt"
Copy-Format "I7" "X74"
$ws.Range("I76").Value = "1024-1125"
Copy-Format "I7" "I76"
$ws.Range("M76").Value = "A Száli-ház"
Copy-Format "M7" "M76"
$ws.Range("I78").Value = "1138-1254"
Copy-Format "I7" "I78"
$ws.Range("M78").Value = "A Stauf-ház"
Copy-Format "M7" "M78"
$ws.Range("X78").Value = "Stauf-ház.txt"
Copy-Format "I7" "X78"
$ws.Range("X76").Value = "Száli-ház.txt"
Copy-Format "I7" "X76"
$ws.Range("I80").Value = "1254-1313"
Copy-Format "I7" "I80"
$ws.Range("M80").Value = "Interregnum és a Habsburgok felemelkedése"
Copy-Format "M7" "M80"
$ws.Range("X80").Value = "Interregnum és Habsburgok.txt"
Copy-Format "I7" "X80"
$ws.Range("I82").Value = "1378-1493"
Copy-Format "I7" "I82"
$ws.Range("M82").Value = "A Habsburg-korszak kezdete"
Copy-Format "M7" "M82"
$ws.Range("X82").Value = "Habsburg.txt"
Copy-Format "I7" "X82"
$ws.Range("I84").Value = "1486-1519"
Copy-Format "I7" "I84"
$ws.Range("M84").Value = "I. Miksa uralkodása"
Copy-Format "M7" "M84"
$ws.Range("X84").Value = "Miksa.txt"
Copy-Format "I7" "X84"
$ws.Range("I86").Value = 1517
Copy-Format "I7" "I86"
$ws.Range("M86").Value = "Reformáció"
Copy-Format "M7" "M86"
$ws.Range("X86").Value = "Reformáció.txt"
Copy-Format "I7" "X86"
$ws.Range("I88").Value = "1618-1648"
Copy-Format "I7" "I88"
$ws.Range("M88").Value = "A harmincéves háború"
Copy-Format "M7" "M88"
$ws.Range("X88").Value = "Harmincéves háború.txt"
Copy-Format "I7" "X88"
$ws.Range("I90").Value = "1799-1815"
Copy-Format "I7" "I90"
$ws.Range("M90").Value = "Napóleoni háborúk"
Copy-Format "M7" "M90"
$ws.Range("X90").Value = "Napóleon.txt"
Copy-Format "I7" "X90"
$ws.Range("I92").Value = "1815-1848"
Copy-Format "I7" "I92"
$ws.Range("M92").Value = "Német szövetség és a Szent Szövetség"
Copy-Format "M7" "M92"
$ws.Range("X92").Value = "Szövetségek.txt"
Copy-Format "I7" "X92"
$ws.Range("I94").Value = 1848
Copy-Format "I7" "I94"
$ws.Range("M94").Value = "48-as német forradalom Berlinben"
Copy-Format "M7" "M94"
$ws.Range("X94").Value = "Német 48-as forradalom"
Copy-Format "I7" "X94"
$ws.Range("I96").Value = 1866
Copy-Format "I7" "I96"
$ws.Range("M96").Value = "porosz-osztrák háború"
Copy-Format "M7" "M96"
$ws.Range("X96").Value = "porosz-osztrák.txt"
Copy-Format "I7" "X96"
$ws.Range("I98").Value = 1871
Copy-Format "I7" "I98"
$ws.Range("M98").Value = "Német Császárság megalakulása"
Copy-Format "M7" "M98"
$ws.Range("I100").Value = 1888
Copy-Format "I7" "I100"
$ws.Range("M100").Value = "Három császár éve"
Copy-Format "M7" "M100"
$ws.Range("X98").Value = "Német Császárság.txt"
Copy-Format "I7" "X98"
$ws.Range("X100").Value = "Három császár éve.txt"
Copy-Format "I7" "X100"

# Row heights for newly added data rows (match style of existing data rows, ht=26.25)
$ws.Rows.Item(72).RowHeight = 26.25
$ws.Rows.Item(74).RowHeight = 26.25
$ws.Rows.Item(76).RowHeight = 26.25
$ws.Rows.Item(78).RowHeight = 26.25
$ws.Rows.Item(80).RowHeight = 26.25
$ws.Rows.Item(82).RowHeight = 26.25
$ws.Rows.Item(84).RowHeight = 26.25
$ws.Rows.Item(86).RowHeight = 26.25
$ws.Rows.Item(88).RowHeight = 26.25
$ws.Rows.Item(90).RowHeight = 26.25
$ws.Rows.Item(92).RowHeight = 26.25
$ws.Rows.Item(94).RowHeight = 26.25
$ws.Rows.Item(96).RowHeight = 26.25
$ws.Rows.Item(98).RowHeight = 26.25
$ws.Rows.Item(100).RowHeight = 26.25

# Update view: selection and scroll position to match target state
$excel.ActiveWindow.ScrollRow = 86
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("X101").Select()
